# Insert a new weekly data row after the existing first data row (row 2),
# pushing the previous rows 3-24 down to 4-25, then populate the new
# row 3 with the new week's values. The columns that are constant across
# all rows for this market/category (A,B,C,E,F,G,H,I,N,O,Q,R) are copied
# from the row that used to be row 3 (now row 4).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(3).Insert()

$ws.Range("A3").Value = $ws.Range("A4").Value()
$ws.Range("B3").Value = $ws.Range("B4").Value()
$ws.Range("C3").Value = $ws.Range("C4").Value()
$ws.Range("D3").Value = 44532
$ws.Range("E3").Value = $ws.Range("E4").Value()
$ws.Range("F3").Value = $ws.Range("F4").Value()
$ws.Range("G3").Value = $ws.Range("G4").Value()
$ws.Range("H3").Value = $ws.Range("H4").Value()
$ws.Range("I3").Value = $ws.Range("I4").Value()
$ws.Range("J3").Value = 300
$ws.Range("K3").Value = 1000
$ws.Range("L3").Value = 1200
$ws.Range("M3").Value = 1100
$ws.Range("N3").Value = $ws.Range("N4").Value()
$ws.Range("O3").Value = $ws.Range("O4").Value()
$ws.Range("P3").Value = 1100
$ws.Range("Q3").Value = $ws.Range("Q4").Value()
$ws.Range("R3").Value = $ws.Range("R4").Value()
